$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.01231499237816
$ws.Range("C2").Value = 11.24466912202767
$ws.Range("D2").Value = 4.65009186814443
$ws.Range("F2").Value = 23.93492906009243
$ws.Range("G2").Value = 3.623931445481563
$ws.Range("I2").Value = 21.86898403228411
$ws.Range("L2").Value = 10.9251793847313
$ws.Range("M2").Value = 14.78889066323006
$ws.Range("N2").Value = 17.77724568101167
$ws.Range("O2").Value = 21.32743305924266
$ws.Range("B3").Value = 13.53893684135647
$ws.Range("C3").Value = 11.04163750199346
$ws.Range("D3").Value = 4.609121112013741
$ws.Range("F3").Value = 23.90269593633092
$ws.Range("G3").Value = 3.62592853167873
$ws.Range("I3").Value = 21.96010210896457
$ws.Range("L3").Value = 10.94063392873699
$ws.Range("M3").Value = 14.70313885070994
$ws.Range("N3").Value = 17.82887292831443
$ws.Range("O3").Value = 21.36052316884369
$ws.Range("B4").Value = 13.24146329985049
$ws.Range("C4").Value = 10.9138922758528
$ws.Range("D4").Value = 4.583423132494522
$ws.Range("F4").Value = 23.89009586802666
$ws.Range("G4").Value = 3.627220840514442
$ws.Range("I4").Value = 22.02075455883773
$ws.Range("L4").Value = 10.95178243925946
$ws.Range("M4").Value = 14.65249574585802
$ws.Range("N4").Value = 17.86234676095697
$ws.Range("O4").Value = 21.38633208164012
$ws.Range("B5").Value = 13.11871276680663
$ws.Range("C5").Value = 10.86110050449506
$ws.Range("D5").Value = 4.572819146868122
$ws.Range("F5").Value = 23.88677334960243
$ws.Range("G5").Value = 3.627764137698688
$ws.Range("I5").Value = 22.04665206835199
$ws.Range("L5").Value = 10.95674297229775
$ws.Range("M5").Value = 14.63237964590724
$ws.Range("N5").Value = 17.87643480416695
$ws.Range("O5").Value = 21.39822722882991
$ws.Range("B6").Value = 13.09824364304865
$ws.Range("C6").Value = 10.85229132763192
$ws.Range("D6").Value = 4.571050528697198
$ws.Range("F6").Value = 23.88633117151324
$ws.Range("G6").Value = 3.627855360190155
$ws.Range("I6").Value = 22.05102360646246
$ws.Range("L6").Value = 10.95759188004873
$ws.Range("M6").Value = 14.62907133281382
$ws.Range("N6").Value = 17.87880115073951
$ws.Range("O6").Value = 21.40028552510316
$ws.Range("B7").Value = 13.2398137699947
$ws.Range("C7").Value = 10.91318322675501
$ws.Range("D7").Value = 4.583280651477257
$ws.Range("F7").Value = 23.89004371851257
$ws.Range("G7").Value = 3.627228100040647
$ws.Range("I7").Value = 22.02109904210888
$ws.Range("L7").Value = 10.9518476484322
$ws.Range("M7").Value = 14.65222232116555
$ws.Range("N7").Value = 17.86253494530445
$ws.Range("O7").Value = 21.38648692953315
$ws.Range("B8").Value = 13.85061608749066
$ws.Range("C8").Value = 11.17532756951123
$ws.Range("D8").Value = 4.63608020507241
$ws.Range("F8").Value = 23.9223253795633
$ws.Range("G8").Value = 3.624606355210897
$ws.Range("I8").Value = 21.89942363492563
$ws.Range("L8").Value = 10.93016382810791
$ws.Range("M8").Value = 14.75891645206545
$ws.Range("N8").Value = 17.79467904422699
$ws.Range("O8").Value = 21.3377010904779
$ws.Range("B9").Value = 14.98693337417091
$ws.Range("C9").Value = 11.6632338803253
$ws.Range("D9").Value = 4.735133079820783
$ws.Range("F9").Value = 24.04243939676645
$ws.Range("G9").Value = 3.619987143477926
$ws.Range("I9").Value = 21.69824929588077
$ws.Range("L9").Value = 10.90079931513579
$ws.Range("M9").Value = 14.98333862285168
$ws.Range("N9").Value = 17.67564947152481
$ws.Range("O9").Value = 21.28571962123335
$ws.Range("B10").Value = 15.77558066578407
$ws.Range("C10").Value = 12.00361662239462
$ws.Range("D10").Value = 4.804926461833069
$ws.Range("F10").Value = 24.16486822285162
$ws.Range("G10").Value = 3.61690833103475
$ws.Range("I10").Value = 21.57339342663986
$ws.Range("L10").Value = 10.88723192110464
$ws.Range("M10").Value = 15.15644837525135
$ws.Range("N10").Value = 17.59669350176222
$ws.Range("O10").Value = 21.27428668822062
$ws.Range("B11").Value = 16.12278676718086
$ws.Range("C11").Value = 12.15412634443512
$ws.Range("D11").Value = 4.83597952702323
$ws.Range("F11").Value = 24.22785202701186
$ws.Range("G11").Value = 3.615575377041547
$ws.Range("I11").Value = 21.52160342521709
$ws.Range("L11").Value = 10.88279392290533
$ws.Range("M11").Value = 15.23675546342955
$ws.Range("N11").Value = 17.56260621098872
$ws.Range("O11").Value = 21.27491020993754
$ws.Range("B12").Value = 16.25249441586474
$ws.Range("C12").Value = 12.21046561160018
$ws.Range("D12").Value = 4.847634211714611
$ws.Range("F12").Value = 24.25273650975909
$ws.Range("G12").Value = 3.615080291303727
$ws.Range("I12").Value = 21.50271430482413
$ws.Range("L12").Value = 10.88136218168998
$ws.Range("M12").Value = 15.26736967251527
$ws.Range("N12").Value = 17.54996049414258
$ws.Range("O12").Value = 21.27598394301803
$ws.Range("B13").Value = 16.22463997125553
$ws.Range("C13").Value = 12.1983615920896
$ws.Range("D13").Value = 4.84512887923961
$ws.Range("F13").Value = 24.24733145754735
$ws.Range("G13").Value = 3.615186487331902
$ws.Range("I13").Value = 21.50675023533617
$ws.Range("L13").Value = 10.88165947464839
$ws.Range("M13").Value = 15.2607676259436
$ws.Range("N13").Value = 17.55267231734195
$ws.Range("O13").Value = 21.27571544788687
$ws.Range("B14").Value = 16.13349398799114
$ws.Range("C14").Value = 12.15877472609042
$ws.Range("D14").Value = 4.836940482346741
$ws.Range("F14").Value = 24.22987865679957
$ws.Range("G14").Value = 3.615534452409814
$ws.Range("I14").Value = 21.52003490881005
$ws.Range("L14").Value = 10.88267114885001
$ws.Range("M14").Value = 15.2392701494631
$ws.Range("N14").Value = 17.56156058573303
$ws.Range("O14").Value = 21.2749817613898
$ws.Range("B15").Value = 16.07743065365197
$ws.Range("C15").Value = 12.13444033173825
$ws.Range("D15").Value = 4.831911125838198
$ws.Range("F15").Value = 24.21932249578483
$ws.Range("G15").Value = 3.615748849631303
$ws.Range("I15").Value = 21.52826633790752
$ws.Range("L15").Value = 10.88332321773216
$ws.Range("M15").Value = 15.22612823434567
$ws.Range("N15").Value = 17.56703905634814
$ws.Range("O15").Value = 21.27464143090417
$ws.Range("B16").Value = 15.75264770491539
$ws.Range("C16").Value = 11.99369045575251
$ws.Range("D16").Value = 4.802882688929557
$ws.Range("F16").Value = 24.16089758716121
$ws.Range("G16").Value = 3.616996799178476
$ws.Range("I16").Value = 21.57687900809999
$ws.Range("L16").Value = 10.88755680407922
$ws.Range("M16").Value = 15.15122973837281
$ws.Range("N16").Value = 17.59895794732396
$ws.Range("O16").Value = 21.27436316320091
$ws.Range("B17").Value = 15.5503619631688
$ws.Range("C17").Value = 11.90621159795965
$ws.Range("D17").Value = 4.78489324608243
$ws.Range("F17").Value = 24.12691329556923
$ws.Range("G17").Value = 3.617779659118505
$ws.Range("I17").Value = 21.60798568054846
$ws.Range("L17").Value = 10.89059773107931
$ws.Range("M17").Value = 15.10566666660623
$ws.Range("N17").Value = 17.61900734574878
$ws.Range("O17").Value = 21.27568448995101
$ws.Range("B18").Value = 15.43293080208676
$ws.Range("C18").Value = 11.85549013038797
$ws.Range("D18").Value = 4.7744809074266
$ws.Range("F18").Value = 24.10805360726673
$ws.Range("G18").Value = 3.618236306353838
$ws.Range("I18").Value = 21.62634859942242
$ws.Range("L18").Value = 10.89251002730636
$ws.Range("M18").Value = 15.07960794142536
$ws.Range("N18").Value = 17.63071153933017
$ws.Range("O18").Value = 21.27699274104024
$ws.Range("B19").Value = 15.39298841868785
$ws.Range("C19").Value = 11.83824801761057
$ws.Range("D19").Value = 4.77094438839609
$ws.Range("F19").Value = 24.10178647600542
$ws.Range("G19").Value = 3.618392014243069
$ws.Range("I19").Value = 21.63264681661709
$ws.Range("L19").Value = 10.89318554765192
$ws.Range("M19").Value = 15.07081093376852
$ws.Range("N19").Value = 17.63470399658218
$ws.Range("O19").Value = 21.27752984337857
$ws.Range("B20").Value = 15.57200842463489
$ws.Range("C20").Value = 11.91556612754562
$ws.Range("D20").Value = 4.786815041482349
$ws.Range("F20").Value = 24.13045996029835
$ws.Range("G20").Value = 3.617695663738087
$ws.Range("I20").Value = 21.60462553790541
$ws.Range("L20").Value = 10.8902571283054
$ws.Range("M20").Value = 15.11050176856155
$ws.Range("N20").Value = 17.61685522531717
$ws.Range("O20").Value = 21.27548708708198
$ws.Range("B21").Value = 16.16031467260013
$ws.Range("C21").Value = 12.170420387402
$ws.Range("D21").Value = 4.839348481238712
$ws.Range("F21").Value = 24.23497703196283
$ws.Range("G21").Value = 3.615431984432856
$ws.Range("I21").Value = 21.51611324600716
$ws.Range("L21").Value = 10.88236724672761
$ws.Range("M21").Value = 15.24557911670399
$ws.Range("N21").Value = 17.55894277046612
$ws.Range("O21").Value = 21.27517453283416
$ws.Range("B22").Value = 16.5344374437494
$ws.Range("C22").Value = 12.33314799290633
$ws.Range("D22").Value = 4.873071456315525
$ws.Range("F22").Value = 24.30930332543116
$ws.Range("G22").Value = 3.614008910660378
$ws.Range("I22").Value = 21.46247855278907
$ws.Range("L22").Value = 10.87866085195645
$ws.Range("M22").Value = 15.33503737046014
$ws.Range("N22").Value = 17.52262275813668
$ws.Range("O22").Value = 21.27985211223697
$ws.Range("B23").Value = 16.33574300256733
$ws.Range("C23").Value = 12.24665827295855
$ws.Range("D23").Value = 4.855130142741559
$ws.Range("F23").Value = 24.26908852575636
$ws.Range("G23").Value = 3.614763289573542
$ws.Range("I23").Value = 21.49071805165184
$ws.Range("L23").Value = 10.88050652513797
$ws.Range("M23").Value = 15.28719091501006
$ws.Range("N23").Value = 17.54186776842314
$ws.Range("O23").Value = 21.27690905643875
$ws.Range("B24").Value = 15.56222558425549
$ws.Range("C24").Value = 11.91133827759372
$ws.Range("D24").Value = 4.785946415386453
$ws.Range("F24").Value = 24.12885439912632
$ws.Range("G24").Value = 3.617733617578605
$ws.Range("I24").Value = 21.60614316564648
$ws.Range("L24").Value = 10.89041060385899
$ws.Range("M24").Value = 15.10831539349315
$ws.Range("N24").Value = 17.61782764590731
$ws.Range("O24").Value = 21.27557462405442
$ws.Range("B25").Value = 14.68707265136703
$ws.Range("C25").Value = 11.53427071157525
$ws.Range("D25").Value = 4.70884262490768
$ws.Range("F25").Value = 24.00390284269143
$ws.Range("G25").Value = 3.621181220838607
$ws.Range("I25").Value = 21.74865181565906
$ws.Range("L25").Value = 10.90733580567144
$ws.Range("M25").Value = 14.92110453225241
$ws.Range("N25").Value = 17.70635394228535
$ws.Range("O25").Value = 21.29508969743879
